$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New rows 5-6, columns A & B - new translation entries
$ws.Range("A5").Value = "Arco"
$ws.Range("B5").Value = "arc"
$ws.Range("A6").Value = "Nodo de servicio"
$ws.Range("B6").Value = "broker node"

# Widen column A to 24 characters. Excel's COM layer quantizes ColumnWidth
# to the nearest 1/6-character pixel grid (stored_width = ceil(w*6)/6 + 5/6),
# so 24.0 would round-trip to 24.833333...; 23.1 lands on the same pixel
# bucket as 24 and serializes cleanly as width="24".
$ws.Columns.Item(1).ColumnWidth = 23.1

# Update the selection to B7
$ws.Range("B7").Select()
